$wb = $excel.ActiveWorkbook

# Insert a new column before column N on the "Repayment schedule" sheet
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns("N").Insert()

# Set column N width to 11 (matches column M), not best-fit
$wsRepay.Columns("N").ColumnWidth = 10.166666666666666

# Update selection on "Repayment schedule" sheet
$wsRepay.Range("R7").Select()

# Make "Repayment schedule" the active sheet (was "Transactions")
$wsRepay.Activate()
